$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    top Heading1 title ("Play Ace Ventura for Free - Slot Game
#    Review"). The new paragraph is plain (no heading style), with
#    "Meta description" bold followed by the rest of the sentence
#    in normal formatting.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$null = $titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover the amusing online slot game Ace Ventura and its various bonus features. Play for free and enjoy the immersive experience.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Near the bottom of the document there used to be a duplicate
#    of the title ("Play Ace Ventura for Free - Slot Game Review",
#    bold) immediately followed by the italic image-prompt
#    paragraph. The bold duplicate-title paragraph is removed
#    entirely (paragraph + mark).
# ------------------------------------------------------------------
$dupTitlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $candText = $cand.Range.Text.TrimEnd([char]13)
    if ($candText -eq "Play Ace Ventura for Free - Slot Game Review" -and $i -ne 1) {
        $dupTitlePara = $cand
        break
    }
}
if ($dupTitlePara -ne $null) {
    $dupTitlePara.Range.Delete()
}

# ------------------------------------------------------------------
# 3) The remaining italic paragraph's text (the old meta-description
#    sentence) is replaced with the new AI image-generation prompt,
#    keeping its italic run formatting intact.
# ------------------------------------------------------------------
$imgPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $candText = $cand.Range.Text.TrimEnd([char]13)
    if ($candText -eq "Discover the amusing online slot game Ace Ventura and its various bonus features. Play for free and enjoy the immersive experience.") {
        $imgPara = $cand
        break
    }
}

$newPromptText = "Create a feature image for Ace Ventura that showcases the game's cartoon style and features a happy Maya warrior with glasses. The image should be colorful and eye-catching, with the Maya warrior standing in the foreground with a big smile, wearing traditional warrior headdress, and holding a magnifying glass in hand. The background should feature elements from the game, such as Ace Ventura characters, animals, or symbols. The image should convey the fun and excitement of the game and appeal to players who enjoy playful and adventurous slot games."

if ($imgPara -ne $null) {
    $target = $d.Range($imgPara.Range.Start, $imgPara.Range.End - 1)
    $target.Text = $newPromptText
}
